# fix bug exeded requeste in google drive
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Bump the price-list date by one day
$ws.Range("A1").Value = 45311

# Correct the unit prices that were wrongly duplicated/inflated
$ws.Range("D29").Value = 169
$ws.Range("D30").Value = 167
$ws.Range("D31").Value = 167
$ws.Range("D32").Value = 167
$ws.Range("D33").Value = 508
